$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '98.150.93'
Set-TextValue $ws.Range('E2') '  +3.07%  '
Set-TextValue $ws.Range('D3') '3.606.72'
Set-TextValue $ws.Range('E3') '  +1.43%  '
Set-TextValue $ws.Range('E4') '  +0.14%  '
Set-TextValue $ws.Range('D5') '243.57'
Set-TextValue $ws.Range('E5') '  +3.68%  '
Set-TextValue $ws.Range('D6') '658.32'
Set-TextValue $ws.Range('E6') '  +1.43%  '
Set-TextValue $ws.Range('D7') '1.71'
Set-TextValue $ws.Range('E7') '  +18.28%  '
Set-TextValue $ws.Range('D8') '0.419'
Set-TextValue $ws.Range('E8') '  +5.74%  '
Set-TextValue $ws.Range('D9') '1.06'
Set-TextValue $ws.Range('E9') '  +7.61%  '
Set-TextValue $ws.Range('D10') '1.00'
Set-TextValue $ws.Range('E10') '  -0.05%  '
Set-TextValue $ws.Range('D11') '3.609.04'
Set-TextValue $ws.Range('E11') '  +1.60%  '
Set-TextValue $ws.Range('D12') '44.19'
Set-TextValue $ws.Range('E12') '  +5.30%  '
Set-TextValue $ws.Range('D13') '0.205'
Set-TextValue $ws.Range('E13') '  +2.30%  '
Set-TextValue $ws.Range('D14') '6.49'
Set-TextValue $ws.Range('E14') '  +0.00%  '
Set-TextValue $ws.Range('B15') 'WrappedBTC'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D15') '97.854.22'
Set-TextValue $ws.Range('E15') '  +3.33%  '
Set-TextValue $ws.Range('B16') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D16') '4.272.59'
Set-TextValue $ws.Range('E16') '  +1.08%  '
Set-TextValue $ws.Range('D17') '0.0000260'
Set-TextValue $ws.Range('E17') '  +3.18%  '
Set-TextValue $ws.Range('B18') 'Polkadot'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D18') '8.66'
Set-TextValue $ws.Range('E18') '  +9.94%  '
Set-TextValue $ws.Range('B19') 'WrappedEther'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D19') '3.617.28'
Set-TextValue $ws.Range('E19') '  +1.59%  '
Set-TextValue $ws.Range('D20') '12.79'
Set-TextValue $ws.Range('E20') '  +1.55%  '
Set-TextValue $ws.Range('D21') '18.02'
Set-TextValue $ws.Range('E21') '  +1.99%  '
Set-TextValue $ws.Range('D22') '0.524'
Set-TextValue $ws.Range('E22') '  +12.47%  '
Set-TextValue $ws.Range('E23') '  +2.54%  '
Set-TextValue $ws.Range('D24') '515.54'
Set-TextValue $ws.Range('E24') '  +2.57%  '
Set-TextValue $ws.Range('D25') '0.0000207'
Set-TextValue $ws.Range('E25') '  +7.56%  '
Set-TextValue $ws.Range('D26') '6.92'
Set-TextValue $ws.Range('E26') '  +5.83%  '
Set-TextValue $ws.Range('D27') '101.31'
Set-TextValue $ws.Range('E27') '  +7.17%  '
Set-TextValue $ws.Range('D28') '13.01'
Set-TextValue $ws.Range('E28') '  +5.21%  '
Set-TextValue $ws.Range('D29') '3.801.84'
Set-TextValue $ws.Range('E29') '  +1.44%  '
Set-TextValue $ws.Range('D30') '0.157'
Set-TextValue $ws.Range('E30') '  +12.84%  '
Set-TextValue $ws.Range('D31') '3.02'
Set-TextValue $ws.Range('E31') '  +0.28%  '
Set-TextValue $ws.Range('D32') '11.87'
Set-TextValue $ws.Range('E32') '  +5.25%  '
Set-TextValue $ws.Range('D33') '0.998'
Set-TextValue $ws.Range('E33') '  -0.20%  '
Set-TextValue $ws.Range('E34') '  +5.11%  '
Set-TextValue $ws.Range('D35') '0.988'
Set-TextValue $ws.Range('E35') '  -0.98%  '
Set-TextValue $ws.Range('D36') '31.84'
Set-TextValue $ws.Range('E36') '  +0.40%  '
Set-TextValue $ws.Range('D37') '8.89'
Set-TextValue $ws.Range('E37') '  +8.11%  '
Set-TextValue $ws.Range('E38') '  +3.44%  '
Set-TextValue $ws.Range('D39') '616.30'
Set-TextValue $ws.Range('E39') '  +9.77%  '
Set-TextValue $ws.Range('D40') '1.66'
Set-TextValue $ws.Range('E40') '  +8.30%  '
Set-TextValue $ws.Range('D41') '2.00'
Set-TextValue $ws.Range('E41') '  +14.27%  '
Set-TextValue $ws.Range('E42') '  +3.47%  '
Set-TextValue $ws.Range('E43') '  -0.03%  '
Set-TextValue $ws.Range('D44') '0.928'
Set-TextValue $ws.Range('E44') '  +3.78%  '
Set-TextValue $ws.Range('D45') '6.01'
Set-TextValue $ws.Range('E45') '  +7.52%  '
Set-TextValue $ws.Range('D46') '0.0439'
Set-TextValue $ws.Range('E46') '  +7.58%  '
Set-TextValue $ws.Range('E47') '  -0.33%  '
Set-TextValue $ws.Range('E48') '  +1.01%  '
Set-TextValue $ws.Range('D49') '8.57'
Set-TextValue $ws.Range('E49') '  +7.41%  '
Set-TextValue $ws.Range('B50') 'EnergySwap'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D50') '32.95'
Set-TextValue $ws.Range('E50') '  -2.70%  '
Set-TextValue $ws.Range('B51') 'dogwifhat'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D51') '3.29'
Set-TextValue $ws.Range('E51') '  +8.89%  '
